$wb = $excel.ActiveWorkbook

# --- GLOBAL RESULTS sheet ---
$ws1 = $wb.Worksheets.Item("GLOBAL RESULTS")

$ws1.Range("C3").Value = 17.537354531331758
$ws1.Range("C5").Value = -0.7366743395896822
$ws1.Range("C7").Value = 41.716158714420196
$ws1.Range("C9").Value = -18.89750712659457
$ws1.Range("C13").Value = 16.598447268422312
$ws1.Range("C15").Value = -0.9425442266839414
$ws1.Range("C17").Value = 17.63088013402975
$ws1.Range("C19").Value = -24.17857563874324
$ws1.Range("C23").Value = 16.598447268422312
$ws1.Range("C25").Value = -0.9425442266839414
$ws1.Range("C27").Value = 17.63088013402975
$ws1.Range("C29").Value = -24.17857563874324
$ws1.Range("C33").Value = 16.598447268422312
$ws1.Range("C35").Value = -0.9425442266839414
$ws1.Range("C37").Value = 17.63088013402975
$ws1.Range("C39").Value = -24.17857563874324
$ws1.Range("C43").Value = 17.134977127082394
$ws1.Range("C45").Value = -0.6942309309452876
$ws1.Range("C47").Value = 31.394189486310243
$ws1.Range("C49").Value = -17.80872939370766
$ws1.Range("C53").Value = 17.02597814532112
$ws1.Range("C55").Value = -0.7970576970723231
$ws1.Range("C57").Value = 28.598097748171803
$ws1.Range("C59").Value = -20.446488633119543
$ws1.Range("C62").Value = 11.352758242551236
$ws1.Range("C63").Value = 28.598097748171803
$ws1.Range("C64").Value = 35.24673678363601
$ws1.Range("C69").Value = 56879.30324254191
$ws1.Range("C70").Value = 2969297.017437632
$ws1.Range("C71").Value = 2912417.7141950903
$ws1.Range("C76").Value = -6560.580896219946

# --- LANDING GEARS sheet ---
$ws9 = $wb.Worksheets.Item("LANDING GEARS")

$ws9.Range("C5").Value = 12.823659727039871
$ws9.Range("C6").Value = 12.823659727039827
$ws9.Range("C7").Value = 16.411049221957363
$ws9.Range("C8").Value = 16.41104922195736
$ws9.Range("C9").Value = 16.41104922195736
$ws9.Range("C10").Value = 16.411049221957356
$ws9.Range("C23").Value = 16.411049221957363
